# Update the 2025Q3 row (row 29) metrics in the recurrence table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C29").Value = 350
$ws.Range("D29").Value = 42
$ws.Range("E29").Value = 308
$ws.Range("F29").Value = 7.228915662650602
